$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold "Play Dr. Jekyll & Mr. Hyde Free Slot by Betsoft -
#    Review 2021" paragraph right before the final "Prompt: ..." paragraph.
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($lastIdx)
$null = $newPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dr. Jekyll &amp; Mr. Hyde Free Slot by Betsoft - Review 2021</w:t></w:r></w:p>")

# 3. Replace the old "Prompt: ..." text with the meta-description copy
#    (minus the "Meta description: " label).
$oldText = 'Prompt: Create a feature image for "Dr. Jekyll & Mr. Hyde" that captures the essence of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses. Use bold and bright colors to make the warrior stand out against the dark laboratory background. Make sure to include symbols from the game, such as the Dr. Jekyll and Mr. Hyde logo, potions, and the laboratory equipment. The warrior should be holding a symbol from the game, such as the potion bottles, key in the lock, or the notebook. Add a touch of humor to the image to reflect the game''s playful yet creepy nature. Make sure the image is eye-catching and appealing to slot game enthusiasts.'
$newText = 'Read our review for Dr. Jekyll & Mr. Hyde free slot by Betsoft. Discover its features, bonuses, and RTP. Play it for free or real money.'
$null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
